$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 17860872
$ws.Cells.Item(19, 9).Value = 31251150
$ws.Cells.Item(19, 10).Value = 7168
$ws.Cells.Item(19, 11).Value = 31251150
$ws.Cells.Item(19, 12).Value = 7168
$ws.Cells.Item(19, 13).Value = -31250975
$ws.Cells.Item(19, 14).Value = -7518

$ws.Cells.Item(113, 8).Value = 6668510.5
$ws.Cells.Item(113, 10).Value = 2200
$ws.Cells.Item(113, 12).Value = 2200
$ws.Cells.Item(113, 14).Value = -8708

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(37, 8).Value = 11342.4
$ws.Cells.Item(37, 10).Value = 15817
$ws.Cells.Item(37, 12).Value = 15817
$ws.Cells.Item(37, 14).Value = -16363

$ws.Cells.Item(55, 8).Value = 14666.444
$ws.Cells.Item(55, 10).Value = 14666.444
$ws.Cells.Item(55, 12).Value = 14666.444
$ws.Cells.Item(55, 14).Value = -15296.444

$ws.Cells.Item(61, 8).Value = 2284746
$ws.Cells.Item(61, 9).Value = 1097230
$ws.Cells.Item(61, 10).Value = 9805680
$ws.Cells.Item(61, 11).Value = 1097230
$ws.Cells.Item(61, 12).Value = 9805680
$ws.Cells.Item(61, 13).Value = -1097018
$ws.Cells.Item(61, 14).Value = -9806104

$ws.Cells.Item(74, 8).Value = 30667896
$ws.Cells.Item(74, 9).Value = 30303682
$ws.Cells.Item(74, 10).Value = 31374904
$ws.Cells.Item(74, 11).Value = 30303682
$ws.Cells.Item(74, 12).Value = 31374904
$ws.Cells.Item(74, 13).Value = -30302808
$ws.Cells.Item(74, 14).Value = -31376652

$ws.Cells.Item(77, 8).Value = 30667896
$ws.Cells.Item(77, 9).Value = 30303682
$ws.Cells.Item(77, 10).Value = 31374904
$ws.Cells.Item(77, 11).Value = 151518410
$ws.Cells.Item(77, 12).Value = 156874520
$ws.Cells.Item(77, 13).Value = -151514042
$ws.Cells.Item(77, 14).Value = -156883256

$ws.Cells.Item(122, 8).Value = 1221.0526
$ws.Cells.Item(122, 9).Value = 1091.2307
$ws.Cells.Item(122, 10).Value = 1502.3334
$ws.Cells.Item(122, 11).Value = 3273.6921
$ws.Cells.Item(122, 12).Value = 4507.0002
$ws.Cells.Item(122, 13).Value = -823.6921000000002
$ws.Cells.Item(122, 14).Value = -9407.0002

$ws.Cells.Item(136, 8).Value = 2284746
$ws.Cells.Item(136, 9).Value = 1097230
$ws.Cells.Item(136, 10).Value = 9805680
$ws.Cells.Item(136, 11).Value = 3291690
$ws.Cells.Item(136, 12).Value = 29417040
$ws.Cells.Item(136, 13).Value = -3289140
$ws.Cells.Item(136, 14).Value = -29422140

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 17228.4
$ws.Cells.Item(82, 9).Value = 1715.3334
$ws.Cells.Item(82, 10).Value = 27570.445
$ws.Cells.Item(82, 11).Value = 1715.3334
$ws.Cells.Item(82, 12).Value = 27570.445
$ws.Cells.Item(82, 13).Value = -1332.3334
$ws.Cells.Item(82, 14).Value = -28336.445

$ws.Cells.Item(85, 8).Value = 17228.4
$ws.Cells.Item(85, 9).Value = 1715.3334
$ws.Cells.Item(85, 10).Value = 27570.445
$ws.Cells.Item(85, 11).Value = 1715.3334
$ws.Cells.Item(85, 12).Value = 27570.445
$ws.Cells.Item(85, 13).Value = -389.3334
$ws.Cells.Item(85, 14).Value = -30222.445

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(51, 8).Value = 8111.3335
$ws.Cells.Item(51, 9).Value = 3000
$ws.Cells.Item(51, 11).Value = 3000
$ws.Cells.Item(51, 13).Value = -2264

$ws.Cells.Item(58, 8).Value = 2066856.5
$ws.Cells.Item(58, 9).Value = 635
$ws.Cells.Item(58, 10).Value = 9092010
$ws.Cells.Item(58, 11).Value = 635
$ws.Cells.Item(58, 12).Value = 9092010
$ws.Cells.Item(58, 13).Value = -432
$ws.Cells.Item(58, 14).Value = -9092416

$ws.Cells.Item(61, 8).Value = 8111.3335
$ws.Cells.Item(61, 9).Value = 3000
$ws.Cells.Item(61, 11).Value = 3000
$ws.Cells.Item(61, 13).Value = -2652

$ws.Cells.Item(68, 8).Value = 17953.908
$ws.Cells.Item(68, 10).Value = 18699.3
$ws.Cells.Item(68, 12).Value = 18699.3
$ws.Cells.Item(68, 14).Value = -20197.3

$ws.Cells.Item(71, 8).Value = 17953.908
$ws.Cells.Item(71, 10).Value = 18699.3
$ws.Cells.Item(71, 12).Value = 56097.89999999999
$ws.Cells.Item(71, 14).Value = -63585.89999999999

$ws.Cells.Item(74, 8).Value = 16146.546
$ws.Cells.Item(74, 10).Value = 16761.2
$ws.Cells.Item(74, 12).Value = 16761.2
$ws.Cells.Item(74, 14).Value = -18509.2

$ws.Cells.Item(77, 8).Value = 16146.546
$ws.Cells.Item(77, 10).Value = 16761.2
$ws.Cells.Item(77, 12).Value = 50283.60000000001
$ws.Cells.Item(77, 14).Value = -59019.60000000001

$ws.Cells.Item(99, 8).Value = 13987.059
$ws.Cells.Item(99, 9).Value = 10897.5
$ws.Cells.Item(99, 10).Value = 16733.334
$ws.Cells.Item(99, 11).Value = 10897.5
$ws.Cells.Item(99, 12).Value = 16733.334
$ws.Cells.Item(99, 13).Value = -9399.5
$ws.Cells.Item(99, 14).Value = -19729.334

$ws.Cells.Item(114, 8).Value = 15342
$ws.Cells.Item(114, 10).Value = 15342
$ws.Cells.Item(114, 12).Value = 15342
$ws.Cells.Item(114, 14).Value = -24020

$ws.Cells.Item(122, 8).Value = 7014.5
$ws.Cells.Item(122, 9).Value = 8256
$ws.Cells.Item(122, 10).Value = 807
$ws.Cells.Item(122, 11).Value = 24768
$ws.Cells.Item(122, 12).Value = 2421
$ws.Cells.Item(122, 13).Value = -22318
$ws.Cells.Item(122, 14).Value = -7321

$ws.Cells.Item(126, 8).Value = 13987.059
$ws.Cells.Item(126, 9).Value = 10897.5
$ws.Cells.Item(126, 10).Value = 16733.334
$ws.Cells.Item(126, 11).Value = 32692.5
$ws.Cells.Item(126, 12).Value = 50200.00199999999
$ws.Cells.Item(126, 13).Value = -30222.5
$ws.Cells.Item(126, 14).Value = -55140.00199999999

$ws.Cells.Item(136, 8).Value = 2066856.5
$ws.Cells.Item(136, 9).Value = 635
$ws.Cells.Item(136, 10).Value = 9092010
$ws.Cells.Item(136, 11).Value = 1905
$ws.Cells.Item(136, 12).Value = 27276030
$ws.Cells.Item(136, 13).Value = 645
$ws.Cells.Item(136, 14).Value = -27281130

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 14).ClearContents()

$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 14).ClearContents()

$ws.Cells.Item(80, 8).Value = 11543.917
$ws.Cells.Item(80, 9).Value = 4665.385
$ws.Cells.Item(80, 10).Value = 19673.092
$ws.Cells.Item(80, 11).Value = 4665.385
$ws.Cells.Item(80, 12).Value = 19673.092
$ws.Cells.Item(80, 13).Value = -3667.385
$ws.Cells.Item(80, 14).Value = -21669.092

$ws.Cells.Item(83, 8).Value = 11543.917
$ws.Cells.Item(83, 9).Value = 4665.385
$ws.Cells.Item(83, 10).Value = 19673.092
$ws.Cells.Item(83, 11).Value = 23326.925
$ws.Cells.Item(83, 12).Value = 98365.46000000001
$ws.Cells.Item(83, 13).Value = -18334.925
$ws.Cells.Item(83, 14).Value = -108349.46

$ws.Cells.Item(102, 8).Value = 8258.546
$ws.Cells.Item(102, 9).Value = 9649.444
$ws.Cells.Item(102, 11).Value = 9649.444
$ws.Cells.Item(102, 13).Value = -8027.444

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2224.8333
$ws.Cells.Item(7, 9).Value = 2386
$ws.Cells.Item(7, 10).Value = 1902.5
$ws.Cells.Item(7, 11).Value = 2386
$ws.Cells.Item(7, 12).Value = 1902.5
$ws.Cells.Item(7, 13).Value = -2274
$ws.Cells.Item(7, 14).Value = -2126.5

$ws.Cells.Item(82, 8).Value = 6013.6665
$ws.Cells.Item(82, 9).Value = 1136
$ws.Cells.Item(82, 10).Value = 8452.5
$ws.Cells.Item(82, 11).Value = 1136
$ws.Cells.Item(82, 12).Value = 8452.5
$ws.Cells.Item(82, 13).Value = -775
$ws.Cells.Item(82, 14).Value = -9174.5

$ws.Cells.Item(85, 8).Value = 6013.6665
$ws.Cells.Item(85, 9).Value = 1136
$ws.Cells.Item(85, 10).Value = 8452.5
$ws.Cells.Item(85, 11).Value = 1136
$ws.Cells.Item(85, 12).Value = 8452.5
$ws.Cells.Item(85, 13).Value = 112
$ws.Cells.Item(85, 14).Value = -10948.5

$ws.Cells.Item(126, 8).Value = 2224.8333
$ws.Cells.Item(126, 9).Value = 2386
$ws.Cells.Item(126, 10).Value = 1902.5
$ws.Cells.Item(126, 11).Value = 7158
$ws.Cells.Item(126, 12).Value = 5707.5
$ws.Cells.Item(126, 13).Value = -4688
$ws.Cells.Item(126, 14).Value = -10647.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(27, 8).Value = 34000
$ws.Cells.Item(27, 10).Value = 34000
$ws.Cells.Item(27, 12).Value = 34000
$ws.Cells.Item(27, 14).Value = -34138

$ws.Cells.Item(115, 8).Value = 50377
$ws.Cells.Item(115, 10).Value = 50377
$ws.Cells.Item(115, 12).Value = 50377
$ws.Cells.Item(115, 14).Value = -53511

$ws.Cells.Item(122, 8).Value = 1808
$ws.Cells.Item(122, 9).Value = 1587.8334
$ws.Cells.Item(122, 10).Value = 2204.3
$ws.Cells.Item(122, 11).Value = 4763.5002
$ws.Cells.Item(122, 12).Value = 6612.900000000001
$ws.Cells.Item(122, 13).Value = -2313.5002
$ws.Cells.Item(122, 14).Value = -11512.9

$ws.Cells.Item(126, 8).Value = 19232730
$ws.Cells.Item(126, 9).Value = 25000750
$ws.Cells.Item(126, 10).Value = 6000
$ws.Cells.Item(126, 11).Value = 75002250
$ws.Cells.Item(126, 12).Value = 18000
$ws.Cells.Item(126, 13).Value = -74999780
$ws.Cells.Item(126, 14).Value = -22940
